$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text runs (Volume number, report week dates) ---
$volChars = $ws.Range("A8").Characters(21, 1)
$volChars.Text = "10"

$weekChars1 = $ws.Range("C9").Characters(27, 9)
$weekChars1.Text = "3/6/2023"
$weekChars2 = $ws.Range("C9").Characters(46, 8)
$weekChars2.Text = "3/12/2023"

# --- Update crime-statistics grid (rows 16-30) ---
$ws.Range("D16").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 15
$ws.Range("K16").Value = -21.052631578947
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -16.666666666666
$ws.Range("N16").Value = -78.873239436619
$ws.Range("C17").Value = 2
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D17").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 13
$ws.Range("K17").Value = -18.75
$ws.Range("L17").Value = 44.444444444444
$ws.Range("M17").Value = 18.181818181818
$ws.Range("N17").Value = -43.478260869565
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 40
$ws.Range("F18").Value = 33
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = 26.923076923076
$ws.Range("I18").Value = 78
$ws.Range("J18").Value = 65
$ws.Range("K18").Value = 20
$ws.Range("L18").Value = 44.444444444444
$ws.Range("M18").Value = 36.842105263157
$ws.Range("N18").Value = -61.386138613861
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 57
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 58.333333333333
$ws.Range("I19").Value = 153
$ws.Range("J19").Value = 99
$ws.Range("K19").Value = 54.545454545454
$ws.Range("L19").Value = 200
$ws.Range("M19").Value = 101.315789473684
$ws.Range("N19").Value = 54.545454545454
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 33.333333333333
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -9.090909090909
$ws.Range("I20").Value = 35
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = 52.173913043478
$ws.Range("L20").Value = 118.75
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = -94.318181818181
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -8.333333333333
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 86
$ws.Range("H21").Value = 31.395348837209
$ws.Range("I21").Value = 295
$ws.Range("J21").Value = 222
$ws.Range("K21").Value = 32.882882882882
$ws.Range("L21").Value = 121.804511278195
$ws.Range("M21").Value = 56.914893617021
$ws.Range("N21").Value = -70.878578479763
$ws.Range("C24").Value = 6
$ws.Range("E24").Value = -53.846153846153
$ws.Range("F24").Value = 39
$ws.Range("H24").Value = -26.415094339622
$ws.Range("I24").Value = 106
$ws.Range("J24").Value = 155
$ws.Range("K24").Value = -31.612903225806
$ws.Range("L24").Value = -5.357142857142
$ws.Range("M24").Value = 6
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 6
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = -25
$ws.Range("I25").Value = 36
$ws.Range("J25").Value = 46
$ws.Range("K25").Value = -21.739130434782
$ws.Range("L25").Value = 176.923076923077
$ws.Range("M25").Value = 33.333333333333
$ws.Range("C27").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -33.333333333333
$ws.Range("L30").Value = -100
$ws.Range("K30").Copy()
$ws.Range("L30").PasteSpecial(-4122)

$excel.CutCopyMode = $false
